$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.649233777652959
$ws.Range("J4").Value = 0.4583536105638675
$ws.Range("K4").Value = 0.3873467178024197
$ws.Range("L4").Value = 2.782278801375452
